$d = $word.ActiveDocument

# Locate the anchor paragraph: "Seria bueno tambien ajustar el tema de las
# imagenes par que sean en otra hoja" -- the new notes get appended right
# after it, before the trailing blank paragraphs.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Seria bueno tambien ajustar el tema de las imágenes par que sean en otra hoja")
$anchor = $searchRange.Paragraphs(1)

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Insert-ParaXml($afterPara, [string]$innerXml) {
    $afterPara.Range.InsertParagraphAfter()
    $newPara = $afterPara.Next()
    $r = $d.Range($newPara.Range.Start, $newPara.Range.End)
    $pkg = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wns + '><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
    return $newPara
}

# Paragraph 1: "Primero verificar como estamanejando lo de las sesiones"
$p1xml = '<w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Primero verificar como </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>estamanejando</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> lo de las sesiones </w:t></w:r></w:p>'
$afterPara = Insert-ParaXml $anchor $p1xml

# Paragraph 2: "Lo primero colocare para que pida el tecnico..."
$p2xml = '<w:p><w:r><w:t xml:space="preserve">Lo primero colocare para que pida el técnico al momento de crear un diagnostico </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>asi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> pues </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>garaztizo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> que lo coloquen de forma correcta </w:t></w:r></w:p>'
$afterPara = Insert-ParaXml $afterPara $p2xml

# Paragraph 3: blank paragraph
$p3xml = '<w:p/>'
$afterPara = Insert-ParaXml $afterPara $p3xml

# Paragraph 4: "Revisar porque cuando voy a grabar una inspeccion..."
$p4xml = '<w:p><w:r><w:t xml:space="preserve">Revisar porque cuando voy a grabar una inspección contraincendios me avisa que falta un campo pero luego ya pasa y lo graba eso no puede ser </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>asi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$afterPara = Insert-ParaXml $afterPara $p4xml

Write-Host "Inserted 4 paragraphs after the 'Seria bueno' note."
